$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 8-15 (columns C, D, E) and relabel column B ---
# Row 8 (was line7 mislabeled extr1 -> now correctly "line7")
$ws.Cells.Item(8,2).Value = "line7"
$ws.Cells.Item(8,3).Value = 14
$ws.Cells.Item(8,4).Value = 11
$ws.Cells.Item(8,5).Value = $true

# Row 9 (now "line8")
$ws.Cells.Item(9,2).Value = "line8"
$ws.Cells.Item(9,3).Value = 16
$ws.Cells.Item(9,4).Value = 9
$ws.Cells.Item(9,5).Value = $true

# Row 10 (now "extr1")
$ws.Cells.Item(10,2).Value = "extr1"
$ws.Cells.Item(10,3).Value = 5
$ws.Cells.Item(10,4).Value = 12
$ws.Cells.Item(10,5).Value = $true

# Row 11 (now "extr2")
$ws.Cells.Item(11,2).Value = "extr2"
$ws.Cells.Item(11,3).Value = 5
$ws.Cells.Item(11,4).Value = 9
$ws.Cells.Item(11,5).Value = $true

# Row 12 (now "extr3")
$ws.Cells.Item(12,2).Value = "extr3"
$ws.Cells.Item(12,3).Value = 10
$ws.Cells.Item(12,4).Value = 11
$ws.Cells.Item(12,5).Value = $false

# Row 13 (now "extr4")
$ws.Cells.Item(13,2).Value = "extr4"
$ws.Cells.Item(13,3).Value = 7
$ws.Cells.Item(13,4).Value = 8
$ws.Cells.Item(13,5).Value = $false

# Row 14 (now "extr5")
$ws.Cells.Item(14,2).Value = "extr5"
$ws.Cells.Item(14,3).Value = 9
$ws.Cells.Item(14,4).Value = 11
$ws.Cells.Item(14,5).Value = $true

# Row 15 (now "extr6")
$ws.Cells.Item(15,2).Value = "extr6"
$ws.Cells.Item(15,3).Value = 7
$ws.Cells.Item(15,4).Value = 11
$ws.Cells.Item(15,5).Value = $true

# --- Add new rows 16 and 17 ---
# Row 16 ("extr7")
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "extr7"
$ws.Cells.Item(16,3).Value = 5
$ws.Cells.Item(16,4).Value = 7
$ws.Cells.Item(16,5).Value = $true

# Row 17 ("extr8")
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "extr8"
$ws.Cells.Item(17,3).Value = 8
$ws.Cells.Item(17,4).Value = 5
$ws.Cells.Item(17,5).Value = $true

# Apply the bold/bordered "index column" style (used by A2:A15) to the new A16/A17 cells
$ws.Cells.Item(15,1).Copy() | Out-Null
$ws.Cells.Item(16,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(15,1).Copy() | Out-Null
$ws.Cells.Item(17,1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
